$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - column F updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1202
$ws1.Range("F3").Value = 602
$ws1.Range("F4").Value = 40
$ws1.Range("F5").Value = 38
$ws1.Range("F6").Value = 208
$ws1.Range("F7").Value = 66
$ws1.Range("F10").Value = 5609
$ws1.Range("F11").Value = 4974
$ws1.Range("F16").Value = 205
$ws1.Range("F17").Value = 9

# Sheet "全部类型" (sheet4) - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1202
$ws4.Range("F3").Value = 602
$ws4.Range("F4").Value = 40
$ws4.Range("F5").Value = 38
$ws4.Range("F6").Value = 208
$ws4.Range("F7").Value = 66
$ws4.Range("F10").Value = 5609
$ws4.Range("F11").Value = 4974
$ws4.Range("F16").Value = 205
$ws4.Range("F19").Value = 9
